# Scene.xlsx - support scene event portal, punish and reward support multi times
#
# Updates the "Quest" (D) column text for several scenes (adding/introducing
# "portal" events alongside the existing sewer/river/fortune/... events) and
# flags new reward/punish columns (E-J) with a 1 where a scene now supports
# that event. Also widens column D to fit the longer text and moves the
# active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D column (Quest) text updates -----------------------------------
$ws.Range("D4").Value  = "sewer;3|river;2|fortune;1|oldtree;1|poppyfield;1"
$ws.Range("D5").Value  = "wolfnest;2|gamble;1|fishpool;2|sewer;3|river;2"
$ws.Range("D6").Value  = "sandpile;1|stone;2"
$ws.Range("D8").Value  = "river;2|stone;3"
$ws.Range("D11").Value = "portal;1|fishpool;1|grave;2"
$ws.Range("D12").Value = "poppyfield;1"
$ws.Range("D13").Value = "poppyfield;1"
$ws.Range("D15").Value = "fortune;1"
$ws.Range("D16").Value = "trees;3|grave;1|portal;1|oldtree;1"
$ws.Range("D18").Value = "river;2|fishpool;1"
$ws.Range("D19").Value = "portal;3"
$ws.Range("D21").Value = "mushroom;1"
$ws.Range("D22").Value = "trees;3|mushroom;1"

# --- E:J flag columns (QPortal/QCardChange/QPiece/QMerchant/QDoctor/QAngel) --
$ws.Range("I5").Value  = 1

$ws.Range("H8").Value  = 1

$ws.Range("F9").Value  = 1

$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1

$ws.Range("J12").Value = 1

$ws.Range("F13").Value = 1
$ws.Range("H13").Value = 1

$ws.Range("J14").Value = 1

$ws.Range("G15").Value = 1

$ws.Range("I16").Value = 1

$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1

$ws.Range("H19").Value = 1

$ws.Range("I21").Value = 1

$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("J22").Value = 1

# --- Column D width (widened to fit the new, longer Quest strings) ---
# NOTE: the host's ColumnWidth setter quantizes to a 1/7-character grid, so
# the nearest input that lands on the target stored width (52.875 "raw" chars)
# is used here.
$ws.Columns.Item(4).ColumnWidth = 52.142857

# --- Selection moves from D17 to D7 -----------------------------------
$ws.Range("D7").Select() | Out-Null
